# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps that get refreshed each time the
# handback status report is (re-)generated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for 47603c2d-....md
$wsOverview.Range("G2").Value = "2016-08-24 15:17:18"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-08-24 15:17:11"
$wsZhCn.Range("K2").Value = "2016-08-24 15:17:29"

# de-de sheet: Correspond Handoff Datetime (shared with Overview!G2 value) /
# Correspond Handback DateTime
$wsDeDe.Range("H2").Value = "2016-08-24 15:17:18"
$wsDeDe.Range("K2").Value = "2016-08-24 15:17:37"
